# Fruta / hortaliza, semanal
# Inserts two new daily-price rows for Kiwi (Macroferia Regional de Talca)
# ahead of the existing data, shifting the rest of the table down by two
# rows (the two rows that fall off the bottom become the new rows 263/264).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 179; everything currently at row 179 and
# below moves down two rows (old 179 -> 181, ..., old 262 -> 264).
$ws.Rows("179:180").Insert()

# --- New row 179 ---
$ws.Cells.Item(179, 1).Value = 5
$ws.Cells.Item(179, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(179, 3).Value = "Maule"
$ws.Cells.Item(179, 4).Value = 44726
$ws.Cells.Item(179, 5).Value = 7
$ws.Cells.Item(179, 6).Value = "Fruta"
$ws.Cells.Item(179, 7).Value = 100101
$ws.Cells.Item(179, 8).Value = "Berries"
$ws.Cells.Item(179, 9).Value = 100101007
$ws.Cells.Item(179, 10).Value = "Kiwi"
$ws.Cells.Item(179, 11).Value = "Hayward"
$ws.Cells.Item(179, 12).Value = "Especial"
$ws.Cells.Item(179, 13).Value = 100
$ws.Cells.Item(179, 14).Value = 10000
$ws.Cells.Item(179, 15).Value = 10000
$ws.Cells.Item(179, 16).Value = 10000
$ws.Cells.Item(179, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(179, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(179, 19).Value = 556
$ws.Cells.Item(179, 20).Value = 18

# --- New row 180 ---
$ws.Cells.Item(180, 1).Value = 5
$ws.Cells.Item(180, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(180, 3).Value = "Maule"
$ws.Cells.Item(180, 4).Value = 44726
$ws.Cells.Item(180, 5).Value = 7
$ws.Cells.Item(180, 6).Value = "Fruta"
$ws.Cells.Item(180, 7).Value = 100101
$ws.Cells.Item(180, 8).Value = "Berries"
$ws.Cells.Item(180, 9).Value = 100101007
$ws.Cells.Item(180, 10).Value = "Kiwi"
$ws.Cells.Item(180, 11).Value = "Hayward"
$ws.Cells.Item(180, 12).Value = "Primera"
$ws.Cells.Item(180, 13).Value = 200
$ws.Cells.Item(180, 14).Value = 8000
$ws.Cells.Item(180, 15).Value = 8000
$ws.Cells.Item(180, 16).Value = 8000
$ws.Cells.Item(180, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(180, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(180, 19).Value = 444
$ws.Cells.Item(180, 20).Value = 18
